# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# -> the "K" column (column G) values get recomputed/rewritten for the
#    existing rows in the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 2
    3  = 1
    4  = 4
    5  = 1
    6  = 1
    7  = 1
    8  = 3
    9  = 1
    10 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
